$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the Treatment query (row 5 / TreatmentTab) ---
# Add an "AND trt.treatment_id IS NOT NULL" condition to the WHERE clause of
# the SQL stored in B5, matching the updated C3DC phs002431 test case query.
$oldWhere = "std.dbgap_accession = 'phs002431' AND dgn.diagnosis = '8041/3 : Small cell carcinoma, NOS'`nORDER BY `n    trt.treatment_id ASC"
$newWhere = "std.dbgap_accession = 'phs002431' AND dgn.diagnosis = '8041/3 : Small cell carcinoma, NOS' AND trt.treatment_id IS NOT NULL`nORDER BY `n    trt.treatment_id ASC"

$treatmentCell = $ws.Range("B5")
$treatmentQuery = $treatmentCell.Value()
$treatmentCell.Value = $treatmentQuery.Replace($oldWhere, $newWhere)

# Re-apply the cell's formatting explicitly so it keeps its wrapped,
# 12pt presentation after the content refresh.
$treatmentCell.WrapText = $true
$treatmentCell.Font.Size = 12

# --- Update the view/selection to match the saved workbook state ---
$ws.Range("C5").Select()
